# Shift Changed to DEPT name
# Rename the per-employee "shift" values (shift1 / shift2) on the
# EmployeeInfo sheet to department-style names (DEPT1 / DEPT2) so the
# shift column matches the department naming convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeInfo")

# Employees previously on "shift1" / "shift2" (rows 2-20) now belong to DEPT1
$dept1Range = $ws.Range("B2:B20")
$dept1Range.Value = "DEPT1"
$dept1Range.ClearFormats()

# Remaining employees (rows 21-26) now belong to DEPT2
$dept2Range = $ws.Range("B21:B26")
$dept2Range.Value = "DEPT2"
$dept2Range.ClearFormats()
